$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TradesLong")
$ws2 = $wb.Worksheets.Item("TradesShort")

# --- Rename the "On exchange" portfolio/account label from "P1" to "A1" ---
# (Column N holds the "On exchange" value in the current layout, before the
#  PT 70% column is removed.)
$ws1.Range("N2").Value = "A1"
$ws2.Range("N2").Value = "A1"

# --- Remove the obsolete "PT 70%" column (column K) entirely ---
# This shifts SL 80%, CreateDate, On exchange and Note one column to the left
# (L->K, M->L, N->M, O->N).
$ws1.Columns.Item(11).Delete()
$ws2.Columns.Item(11).Delete()

# --- Update the "SL 80%" value (now column K) ---
# TradesLong keeps the previously-computed stop-loss value.
$ws1.Range("K2").Value = 19840
# TradesShort gets a freshly computed stop-loss value (logic fix for short trades).
$ws2.Range("K2").Value = 70060

# --- Selection / active sheet bookkeeping ---
$ws1.Range("K2").Select()
$ws2.Activate()
$ws2.Range("G23").Select()
